$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 14-15 (pushes the "employee_Email" block and everything
# below it down by two rows), matching the diff where new employee_first /
# employee_last attribute rows are added right after the "employee_JobTitle"
# header row (row 13) and before "employee_Email" (old row 14).
$ws.Rows("14:15").Insert()

# Populate the new rows in the exact order needed so that the newly created
# shared-string entries land on the same indices as the target workbook
# (employee_last=140, employee_first=141, "Employee's first name"=142,
#  "Employee's last name"=143).
$ws.Range("B15").Value = "employee_last"
$ws.Range("B14").Value = "employee_first"
$ws.Range("C14").Value = "Employee's first name"
$ws.Range("C15").Value = "Employee's last name"

$ws.Range("D14").Value = "VARCHAR(20)"
$ws.Range("E14").Value = "Xxxxxxxx"
$ws.Range("G14").Value = "No"

$ws.Range("D15").Value = "VARCHAR(20)"
$ws.Range("E15").Value = "Xxxxxxxx"
$ws.Range("G15").Value = "No"

# Match formatting used by the other VARCHAR(20)/Xxxxxxxx rows in the table.
$ws.Range("D34").Copy()
$ws.Range("D14:D15").PasteSpecial(-4122)

$ws.Range("E18").Copy()
$ws.Range("E14:E15").PasteSpecial(-4122)

$ws.Range("B34").Copy()
$ws.Range("B14:B15").PasteSpecial(-4122)
$ws.Range("C14:C15").PasteSpecial(-4122)
$ws.Range("G14:G15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Resize the data table (Table2 / ListObject) and its AutoFilter so the range
# grows from A1:I54 to A1:I56 to include the two new rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:I56"))

# The mailto hyperlink on the "employee_Email" row needs to move from E14 to
# E16 since that row shifted down by two.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E16"), "mailto:xxxx@autoparts.com")

# Restore the worksheet selection recorded in the saved file.
$ws.Range("G18").Select()
